$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (bold, bordered, centered) from A1 onto
# the three new header cells, then set their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Fill in the season record (Wins / Losses / Ties) for every player row.
$ws.Range("AD2:AD56").Value = 79
$ws.Range("AE2:AE56").Value = 83
$ws.Range("AF2:AF56").Value = 0

Write-Host "Season record columns added"
